$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the existing (default) style of the Price column, then force
# it to text so numeric-looking strings (e.g. "592.24") are not silently
# converted into Number cells by Excel's automatic type inference.
$priceRange = $ws.Range("D2:D51")
$origPriceStyle = $ws.Range("D2").Style
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "67.336.19"
$ws.Range("E2").Value = "  -4.61%  "

$ws.Range("D3").Value = "3.255.44"
$ws.Range("E3").Value = "  -7.30%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "592.24"
$ws.Range("E5").Value = "  -4.45%  "

$ws.Range("D6").Value = "150.92"
$ws.Range("E6").Value = "  -12.90%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "3.247.25"
$ws.Range("E8").Value = "  -7.44%  "

$ws.Range("D9").Value = "0.542"
$ws.Range("E9").Value = "  -11.14%  "

$ws.Range("D10").Value = "0.171"
$ws.Range("E10").Value = "  -13.24%  "

$ws.Range("D11").Value = "6.79"
$ws.Range("E11").Value = "  -3.54%  "

$ws.Range("D12").Value = "0.505"
$ws.Range("E12").Value = "  -13.13%  "

$ws.Range("D13").Value = "38.15"
$ws.Range("E13").Value = "  -17.61%  "

$ws.Range("D14").Value = "0.0000243"
$ws.Range("E14").Value = "  -12.10%  "

$ws.Range("D15").Value = "3.774.78"
$ws.Range("E15").Value = "  -7.43%  "

$ws.Range("D16").Value = "67.366.22"
$ws.Range("E16").Value = "  -4.63%  "

$ws.Range("D17").Value = "546.24"
$ws.Range("E17").Value = "  -10.67%  "

$ws.Range("D18").Value = "3.256.13"
$ws.Range("E18").Value = "  -7.55%  "

$ws.Range("D19").Value = "7.25"
$ws.Range("E19").Value = "  -13.58%  "

$ws.Range("E20").Value = "  -6.00%  "

$ws.Range("D21").Value = "15.14"
$ws.Range("E21").Value = "  -14.43%  "

$ws.Range("D22").Value = "0.763"
$ws.Range("E22").Value = "  -13.47%  "

$ws.Range("D23").Value = "7.83"
$ws.Range("E23").Value = "  -13.98%  "

$ws.Range("D24").Value = "85.37"
$ws.Range("E24").Value = "  -13.14%  "

$ws.Range("D25").Value = "13.46"
$ws.Range("E25").Value = "  -13.40%  "

$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("D27").Value = "3.22"
$ws.Range("E27").Value = "  -13.80%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "8.03"
$ws.Range("E28").Value = "  -11.20%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "29.29"
$ws.Range("E29").Value = "  -12.98%  "

$ws.Range("D30").Value = "2.12"
$ws.Range("E30").Value = "  -17.35%  "

$ws.Range("D31").Value = "2.66"
$ws.Range("E31").Value = "  -11.90%  "

$ws.Range("E32").Value = "  -12.37%  "

$ws.Range("D33").Value = "543.70"
$ws.Range("E33").Value = "  -14.92%  "

$ws.Range("D34").Value = "6.64"
$ws.Range("E34").Value = "  -17.74%  "

$ws.Range("D35").Value = "5.69"
$ws.Range("E35").Value = "  -15.88%  "

$ws.Range("E36").Value = "  +0.08%  "

$ws.Range("D37").Value = "0.0449"
$ws.Range("E37").Value = "  -5.35%  "

$ws.Range("D38").Value = "53.61"

$ws.Range("D39").Value = "0.0855"
$ws.Range("E39").Value = "  -14.17%  "

$ws.Range("D40").Value = "9.17"
$ws.Range("E40").Value = "  -14.89%  "

$ws.Range("E41").Value = "  -11.28%  "

$ws.Range("D42").Value = "2.933.07"
$ws.Range("E42").Value = "  -12.23%  "

$ws.Range("D43").Value = "2.61"
$ws.Range("E43").Value = "  -23.36%  "

$ws.Range("D44").Value = "0.261"
$ws.Range("E44").Value = "  -16.21%  "

$ws.Range("D45").Value = "0.0₃0580"
$ws.Range("E45").Value = "  -18.95%  "

$ws.Range("D46").Value = "26.40"
$ws.Range("E46").Value = "  -16.90%  "

$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "128.27"
$ws.Range("E47").Value = "  -4.07%  "

$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").Value = "2.14"
$ws.Range("E48").Value = "  -15.25%  "

$ws.Range("E49").Value = "  -0.08%  "

$ws.Range("D50").Value = "2.35"
$ws.Range("E50").Value = "  -20.34%  "

$ws.Range("E51").Value = "  -12.34%  "

# Restore the original (General) style/number format on the Price column
# now that the text values are safely stored.
$priceRange.Style = $origPriceStyle
